$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45171 -> 45172, i.e. 2023-09-02 -> 2023-09-03) for every data row
# (rows 2 through 303).
$ws.Range("C2:C303").Value = 45172
